# Section 3.2 -> 3.2.2 edit:
#  1. Update the cached "datetimeFigureOut" date field text (shown via the
#     Date placeholder) on the slide master and every slide layout from
#     04/03/2015 to 03/09/2015.
#  2. Give both "Straight Arrow Connector" shapes on slide 1 a triangle
#     arrowhead at the tail end (they already have one at the head end).

$p = $ppt.ActivePresentation

# ppPlaceholderDate
$ppPlaceholderDate = 16
$NewDateText = "03/09/2015"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $shape.TextFrame.TextRange.Text = $NewDateText
            }
        }
    }
}

# Slide master's own Date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date placeholder ("Apply to All" style update).
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Add a tail arrowhead to the two straight arrow connectors on slide 1.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.Name -like "Straight Arrow Connector*") {
        $shape.Line.EndArrowheadStyle = 2   # msoArrowheadTriangle
    }
}
